# Commit message: "add => to URL"
#
# The resume's last bullet currently reads:
#   "  - staging URL http://howardthurmanfilm.herokuapp.com"
# and needs to become:
#   "  - staging URL => http://howardthurmanfilm.herokuapp.com"
#
# i.e. insert the literal characters "=> " right before the "http://" so the
# line reads "staging URL =>" followed by the address.

$d = $word.ActiveDocument

$d.Content.Find.ClearFormatting()
$found = $d.Content.Find.Execute(
    "staging URL http://howardthurmanfilm.herokuapp.com", # FindText
    $true,                                                 # MatchCase
    $false,                                                # MatchWholeWord
    $false,                                                # MatchWildcards
    $false,                                                # MatchSoundsLike
    $false,                                                # MatchAllWordForms
    $true,                                                 # Forward
    1,                                                      # Wrap (wdFindContinue)
    $false,                                                 # Format
    "staging URL => http://howardthurmanfilm.herokuapp.com", # ReplaceWith
    2                                                       # Replace (wdReplaceAll)
)

if (-not $found) {
    throw "Could not locate the staging URL bullet to update."
}

Write-Host "Updated staging URL line with '=>' : $found"
